# Apply cryptos list price/volume updates (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.173.53"
$ws.Range("E2").Value = "  -0.19%  "
# Row 3
$ws.Range("D3").Value = "2.498.74"
$ws.Range("E3").Value = "  -1.32%  "
# Row 4
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$c = $ws.Range("D5")
$c.Value = "'318.26"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "
# Row 6
$c = $ws.Range("D6")
$c.Value = "'105.90"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.09%  "
# Row 7
$ws.Range("E7").Value = "  -1.74%  "
# Row 8
$ws.Range("E8").Value = "  -0.01%  "
# Row 9
$ws.Range("E9").Value = "  -4.06%  "
# Row 10
$c = $ws.Range("D10")
$c.Value = "'38.85"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.78%  "
# Row 11
$c = $ws.Range("D11")
$c.Value = "'20.23"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "
# Row 12
$c = $ws.Range("D12")
$c.Value = "'0.0803"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.40%  "
# Row 13
$ws.Range("E13").Value = "  -0.06%  "
# Row 14
$c = $ws.Range("D14")
$c.Value = "'7.10"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.70%  "
# Row 15
$ws.Range("D15").Value = "2.892.68"
$ws.Range("E15").Value = "  -0.99%  "
# Row 16
$ws.Range("D16").Value = "2.503.07"
$ws.Range("E16").Value = "  -1.29%  "
# Row 17
$c = $ws.Range("D17")
$c.Value = "'0.828"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.88%  "
# Row 18
$ws.Range("D18").Value = "48.026.14"
$ws.Range("E18").Value = "  -0.12%  "
# Row 19
$ws.Range("E19").Value = "  +11.27%  "
# Row 20
$c = $ws.Range("D20")
$c.Value = "'12.86"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.41%  "
# Row 21
$c = $ws.Range("D21")
$c.Value = "'6.59"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
# Row 22
$ws.Range("D22").Value = "0.0₃0930"
$ws.Range("E22").Value = "  -2.70%  "
# Row 23
$c = $ws.Range("D23")
$c.Value = "'71.16"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.78%  "
# Row 24
$c = $ws.Range("D24")
$c.Value = "'267.33"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.39%  "
# Row 25
$c = $ws.Range("D25")
$c.Value = "'2.51"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.62%  "
# Row 26
$ws.Range("E26").Value = "  +0.19%  "
# Row 27
$c = $ws.Range("D27")
$c.Value = "'25.78"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
# Row 28
$c = $ws.Range("D28")
$c.Value = "'2.24"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "
# Row 29
$c = $ws.Range("D29")
$c.Value = "'9.73"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "
# Row 30
$ws.Range("E30").Value = "  -4.04%  "
# Row 31
$c = $ws.Range("D31")
$c.Value = "'34.59"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.23%  "
# Row 32
$c = $ws.Range("D32")
$c.Value = "'49.42"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "
# Row 33
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D33")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
# Row 34
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D34")
$c.Value = "'19.15"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.28%  "
# Row 35
$c = $ws.Range("D35")
$c.Value = "'5.29"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.35%  "
# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.0773"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.79%  "
# Row 37
$c = $ws.Range("D37")
$c.Value = "'1.94"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.84%  "
# Row 38
$c = $ws.Range("D38")
$c.Value = "'4.60"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "
# Row 39
$c = $ws.Range("D39")
$c.Value = "'2.88"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.41%  "
# Row 40
$c = $ws.Range("D40")
$c.Value = "'123.37"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.25%  "
# Row 41
$c = $ws.Range("D41")
$c.Value = "'22.44"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
# Row 42
$c = $ws.Range("D42")
$c.Value = "'0.110"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "
# Row 43
$ws.Range("E43").Value = "  +1.36%  "
# Row 44
$ws.Range("E44").Value = "  +0.59%  "
# Row 45
$ws.Range("D45").Value = "1.999.45"
$ws.Range("E45").Value = "  -0.82%  "
# Row 46
$c = $ws.Range("D46")
$c.Value = "'3.14"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
# Row 47
$ws.Range("E47").Value = "  +1.35%  "
# Row 48
$ws.Range("E48").Value = "  -2.71%  "
# Row 49
$c = $ws.Range("D49")
$c.Value = "'8.95"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.32%  "
# Row 50
$c = $ws.Range("D50")
$c.Value = "'5.20"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "
# Row 51
$c = $ws.Range("D51")
$c.Value = "'78.89"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
